# Update marksheet corrected/total marks values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row - right-answer marking value changed 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row - total right-answer marks changed 66 -> 110
$ws.Range("B12").Value = 110

# "Total" row - correct/total marks string changed 64/84 -> 110/140
$ws.Range("E12").Value = "110/140"
